$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 3.6
$ws.Range("I2").Value = 2.2
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("Q2").Value = 1.85
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 2.4
$ws.Range("T2").Value = 1.53
$ws.Range("U2").Value = 3.85
$ws.Range("AA2").Value = 2
$ws.Range("AB2").Value = 1.75
$ws.Range("AN2").Value = 6.5

# Row 4 updates
$ws.Range("N4").Value = 7.5
$ws.Range("AS4").Value = 67

# Row 5 updates
$ws.Range("L5").Value = 3.6
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 7.5
$ws.Range("AI5").Value = 7.5
$ws.Range("AL5").Value = 51
$ws.Range("AM5").Value = 401

# Row 8 updates
$ws.Range("H8").Value = 3.5
$ws.Range("I8").Value = 2.05
$ws.Range("J8").Value = 3.75
$ws.Range("O8").Value = 1.25
$ws.Range("P8").Value = 3.55
$ws.Range("S8").Value = 1.75
$ws.Range("T8").Value = 1.95
$ws.Range("W8").Value = 2.8
$ws.Range("Y8").Value = 1.38
$ws.Range("Z8").Value = 2.82
$ws.Range("AB8").Value = 2.12
$ws.Range("AG8").Value = 27
$ws.Range("AJ8").Value = 6.8
$ws.Range("AR8").Value = 15.5

# Row 9 updates - fill in previously empty odds cells with new values
$ws.Range("G9").Value = 10.5
$ws.Range("H9").Value = 6.2
$ws.Range("I9").Value = 1.21
$ws.Range("J9").Value = 7.7
$ws.Range("K9").Value = 2.92
$ws.Range("L9").Value = 1.57
$ws.Range("O9").Value = 1.09
$ws.Range("P9").Value = 6
$ws.Range("S9").Value = 1.31
$ws.Range("T9").Value = 3.15
$ws.Range("W9").Value = 1.78
$ws.Range("X9").Value = 1.93
$ws.Range("Y9").Value = 1.2
$ws.Range("Z9").Value = 4.05
$ws.Range("AA9").Value = 1.72
$ws.Range("AB9").Value = 2
$ws.Range("AC9").Value = 40
$ws.Range("AD9").Value = 100
$ws.Range("AE9").Value = 32
$ws.Range("AF9").Value = 300
$ws.Range("AG9").Value = 110
$ws.Range("AH9").Value = 70
$ws.Range("AI9").Value = 28
$ws.Range("AJ9").Value = 14
$ws.Range("AK9").Value = 21
$ws.Range("AL9").Value = 70
$ws.Range("AM9").Value = 350
$ws.Range("AN9").Value = 11.75
$ws.Range("AO9").Value = 8.25
$ws.Range("AP9").Value = 9.75
$ws.Range("AQ9").Value = 8.5
$ws.Range("AR9").Value = 9.75
$ws.Range("AS9").Value = 22
